$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the rows for years 2000, 2005-2009 (rows 2-7), shifting 2010-2013 data up.
$ws.Rows("2:7").Delete()
